# Move list.xlsx update
# 1) Update two "Noca Language" cells with new parser syntax strings.
# 2) Add three new moves (Ember, Vine Whip, Dragon Pride) as new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update existing Noca Language strings ---
# Row 48 (Reflect) column G
$ws.Range("G48").Value = "apply status{reflect} for target in range{Allies:5}"
# Row 56 (Twister) column G
$ws.Range("G56").Value = "apply damage{25} if chance{33%} for target in range{Enemy:5}"

# --- 2) Add three new move rows at the bottom of the table (rows 73-75) ---
# Copy formatting of the last existing data row (row 72) down into the
# three new rows so they inherit the same cell style.
$ws.Range("A72:G72").Copy()
$ws.Range("A73:G73").Insert()
$ws.Range("A72:G72").Copy()
$ws.Range("A74:G74").Insert()
$ws.Range("A72:G72").Copy()
$ws.Range("A75:G75").Insert()

# Row 73: Ember
$ws.Range("A73").Value = "Ember"
$ws.Range("B73").Value = "Lesser Special"
$ws.Range("C73").Value = 0.0
$ws.Range("D73").Value = "2 Tiles"
$ws.Range("E73").Value = "2d6 + 0.5 INT fire damage"
$ws.Range("F73").Value = "Ember"
$ws.Range("G73").Value = "none"

# Row 74: Vine Whip
$ws.Range("A74").Value = "Vine Whip"
$ws.Range("B74").Value = "Lesser Physical"
$ws.Range("C74").Value = 0.0
$ws.Range("D74").Value = "2 Tiles"
$ws.Range("E74").Value = "1d12 + 0.5 INT grass damage"
$ws.Range("F74").Value = "Vine Whip"
$ws.Range("G74").Value = "none"

# Row 75: Dragon Pride
$ws.Range("A75").Value = "Dragon Pride"
$ws.Range("B75").Value = "Lesser Status"
$ws.Range("C75").Value = "Self"
$ws.Range("D75").Value = "Self"
$ws.Range("E75").Value = "None"
$ws.Range("F75").Value = "Dragon Pride"
$ws.Range("G75").Value = "none"
